$d = $word.ActiveDocument

# Remove the " (Rough Draft)" text that was appended after "A: Planning"
$d.Content.Find.Execute(" (Rough Draft)", $false, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
